$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('good child', 'いい子|いいこ')
    ,@('color', '色|いろ')
    ,@('boxed lunch', 'お弁当|おべんとう')
    ,@('Kabuki; traditional Japanese theatrical art', '歌舞伎|かぶき')
    ,@('guitar', 'ギター')
    ,@('last year', '去年|きょねん')
    ,@('medicine', '薬|くすり')
    ,@('to take medicine', '薬を飲む|くすりをのむ')
    ,@('concert', 'コンサート')
    ,@('near future', '今度|こんど')
    ,@('essay; composition', '作文|さくぶん')
    ,@('exam', '試験|しけん')
    ,@('ski', 'スキー')
    ,@('last month', '先月|せんげつ')
    ,@('word; vocabulary', '単語|たんご')
    ,@('piano', 'ピアノ')
    ,@('pizza', 'ピザ')
    ,@('illness; sickness', '病気|びょうき')
    ,@('blue', '青い|あおい')
    ,@('red', '赤い|あかい')
    ,@('black', '黒い|くろい')
    ,@('lonely', '寂しい|さびしい')
    ,@('white', '白い|しろい')
    ,@('young', '若い|わかい')
    ,@('mean-spirited', '意地悪|いじわる(な)')
    ,@('to dance', '踊る|おどる')
    ,@('(something) ends', '終わる|おわる')
    ,@('to be popular', '人気がある|にんきがある')
    ,@('(something) begins', '始まる|はじまる')
    ,@('to play (a string instrument or piano)', '弾く|ひく')
    ,@('to get (from somebody)', 'もらう')
    ,@('to memorize', '覚える|おぼえる')
    ,@('to appear; to attend; to exit', '出る|でる')
    ,@('to do physical exercises', '運動する|うんどうする')
    ,@('to take a walk', '散歩する|さんぽする')
    ,@('from...', '～から')
    ,@('by all means', '是非|ぜひ')
    ,@('by the way', 'ところで')
    ,@('all', 'みんな')
    ,@('already', 'もう')
    ,@('one', '一つ|ひとつ')
    ,@('two', '二つ|ふたつ')
    ,@('three', '三つ|みっつ')
    ,@('four', '四つ|よっつ')
    ,@('five', '五つ|いつつ')
    ,@('six', '六つ|むっつ')
    ,@('seven', '七つ|ななつ')
    ,@('eight', '八つ|やっつ')
    ,@('nine', '九つ|ここのつ')
    ,@('ten', '十|とお')
    ,@('black', '黒い|くろい')
    ,@('white', '白い|しろい')
    ,@('red', '赤い|あかい')
    ,@('blue', '青い|あおい')
    ,@('yellow', '黄色い|きいろい')
    ,@('brown', '茶色い|ちゃいろい')
    ,@('green', '緑／グリーン|みどり／---')
    ,@('purple', '紫|むらさき')
    ,@('gray', '灰色／グレー|はいいろ／---')
    ,@('light blue', '水色|みずいろ')
    ,@('pink', 'ピンク')
    ,@('gold', '金色／ゴールド|きんいろ／---')
    ,@('silver', '銀色／シルバー|ぎんいろ／---')
    ,@('color', '色|いろ')
    ,@('red bag', '赤いかばん')
    ,@('I like red the best.', '赤がいちばん好きです。')
    ,@('black cat', '黒猫|くろねこ')
    ,@('green sweater', '緑／グリーンのセーター')
    ,@('Your face looks pale.', '顔が青いですね。|かおがあおいですね。')
    ,@('black and white picture', '白黒の写真|しろくろのしゃしん')
    ,@('Mary has blonde hair.', 'メアリーさんは金髪です。|メアリーさんはきんぱつです。')
    ,@('A.M.', '午前|ごぜん')
    ,@('P.M.; in the afternoon', '午後|ごご')
    ,@('in the morning', '午前中|ごぜんちゅう')
    ,@('noon', '正午|しょうご')
    ,@('after...', '～の後|～のあと')
    ,@('later', '後で|あとで')
    ,@('back; behind', '後ろ|うしろ')
    ,@('lastly', '最後に|さいごに')
    ,@('before; front', '前|まえ')
    ,@('name', '名前|なまえ')
    ,@('advance sale', '前売り|まえうり')
    ,@('name', '名前|なまえ')
    ,@('famous', '有名な|ゆうめいな')
    ,@('name card', '名刺|めいし')
    ,@('full name', '氏名|しめい')
    ,@('place name', '地名|ちめい')
    ,@('white', '白い|しろい')
    ,@('blank sheet', '白紙|はくし')
    ,@('white color', '白|しろ')
    ,@('swan', '白鳥|はくちょう')
    ,@('rain', '雨|あめ')
    ,@('rainy season', '雨期|うき')
    ,@('rainy season (East Asia)', '梅雨|つゆ')
    ,@('to write', '書く|かく')
    ,@('dictionary', '辞書|じしょ')
    ,@('textbook', '教科書|きょうかしょ')
    ,@('library', '図書館|としょかん')
    ,@('friend', '友だち|ともだち')
    ,@('best friend', '親友|しんゆう')
    ,@('friend (formal)', '友人|ゆうじん')
    ,@('friendship', '友情|ゆうじょう')
    ,@('time', '時間|じかん')
    ,@('two hours', '二時間|にじかん')
    ,@('between', '間|あいだ')
    ,@('human being', '人間|にんげん')
    ,@('one week', '一週間|いっしゅうかん')
    ,@('house', '家|いえ')
    ,@('family', '家族|かぞく')
    ,@('house; home', '家|うち')
    ,@('my wife', '家内|かない')
    ,@('author', '作家|さっか')
    ,@('to speak', '話す|はなす')
    ,@('talk; story', '話|はなし')
    ,@('telephone', '電話|でんわ')
    ,@('conversation', '会話|かいわ')
    ,@('little', '少し|すこし')
    ,@('few', '少ない|すくない')
    ,@('a little', '少々|しょうしょう')
    ,@('girl', '少女|しょうじょ')
    ,@('boy', '少年|しょうねん')
    ,@('old (for things)', '古い|ふるい')
    ,@('secondhand', '中古|ちゅうこ')
    ,@('ancient times', '古代|こだい')
    ,@('to know', '知る|しる')
    ,@('acquaintance (formal)', '知人|ちじん')
    ,@('acquaintance', '知り合い|しりあい')
    ,@('to come', '来る|くる')
    ,@('to come (long-form)', '来ます|きます')
    ,@('not to come', '来ない|こない')
    ,@('next week', '来週|らいしゅう')
    ,@('visit to Japan', '来日|らいにち')
    ,@('Thank you for everything.', 'いろいろおせわになりました。')
    ,@('Please take care of yourself.', '体に気をつけてください。')
    ,@('I am looking forward to seeing you.', 'お会いできるのを楽しみにしています。')
    ,@('Congratulations on...', '～おめでとう（ございます）。')
    ,@('Happy birthday.', '（お）たんじょうびおめでとう。')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}